# Replace the four "2018 Fechas de la campaña para Perseo: Del 30 de ..."
# paragraphs with the translated "Fechas de la campaña para Perseo: Taurus:
# 16-25 de enero" text, per the commit diff.

$d = $word.ActiveDocument
$newText = "Fechas de la campaña para Perseo: Taurus: 16-25 de enero"

$targetIndexes = @(3, 53, 86, 121)

foreach ($idx in $targetIndexes) {
    $p = $d.Paragraphs.Item($idx)
    $r = $p.Range
    # Shrink the range so it excludes the trailing paragraph-mark character,
    # then remove the old run content and insert the new plain-text run.
    $r.MoveEnd(1, -1)
    $r.Delete()

    $p2 = $d.Paragraphs.Item($idx)
    $r2 = $p2.Range
    $r2.InsertAfter($newText)
}
